# Update the hearing/suspension dates from June 10, 2022 -> June 11, 2022,
# and the signature date from August 09, 2022 -> August 10, 2022.

$d = $word.ActiveDocument

# 1) " on June 10, 2022." -> " on June 11, 2022."
$d.Content.Find.Execute(" on June 10, 2022.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " on June 11, 2022.", 2) | Out-Null

# 2) Standalone bold "June 10, 2022" -> "June 11, 2022"
$d.Content.Find.Execute("June 10, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "June 11, 2022", 2) | Out-Null

# 3) "August 09, 2022" -> "August 10, 2022"
$d.Content.Find.Execute("August 09, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "August 10, 2022", 2) | Out-Null

# 4) " license is suspended from June 10, 2022" -> " license is suspended from June 11, 2022"
$d.Content.Find.Execute(" license is suspended from June 10, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " license is suspended from June 11, 2022", 2) | Out-Null
